$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SolverSettings")

# Append the new "MaxLoan_yrs" scenario-parameter row beneath the existing
# SolverSettings table (previously rows 1-8, now rows 1-9).
$ws.Range("A9").Value = "MaxLoan_yrs"
$ws.Range("B9").Value = 30
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 30

# Match the author's resulting selection (one cell past the new row).
[void]$ws.Range("F9").Select()
